# Applies the edits described by the commit:
#  1. Three tables (on slides 14, 15, 16) switch from the custom
#     "Table_0" style ({9233B89D-1A64-452A-A378-F084BC382A85}) to the
#     built-in table style {5248F683-36FA-4675-B9DA-0A894776FB77}.
#  2. The presentation's (slide-master) theme colour palette is swapped
#     from the "Integral" / "Red Violet" scheme to the default
#     "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$newStyleId = "{5248F683-36FA-4675-B9DA-0A894776FB77}"
$tableSlides = @(14, 15, 16)

foreach ($slideIdx in $tableSlides) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Recolour the theme from "Integral"/Red Violet to Office -----------
$officeColors = @(
    0x000000, # dk1
    0xFFFFFF, # lt1
    0x6A5444, # dk2      (44546A)
    0xE6E6E7, # lt2      (E7E6E6)
    0xD59B5B, # accent1  (5B9BD5)
    0x317DED, # accent2  (ED7D31)
    0xA5A5A5, # accent3  (A5A5A5)
    0x00C0FF, # accent4  (FFC000)
    0xC47244, # accent5  (4472C4)
    0x47AD70, # accent6  (70AD47)
    0xC16305, # hlink    (0563C1)
    0x724F95  # folHlink (954F72)
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
